{"js": "// The \"Due:\" date is being updated from \"July 13, 2022\" to \"March 21, 2023\".\n// (All other changes in the source diff are auto-generated proofing marks\n// (w:proofErr spell/grammar check tags) and OOXML namespace/schema bumps\n// produced by a newer Word version opening/re-saving the file; they carry\n// no visible text change, so there is nothing else to edit here.)\n\nconst searchResults = context.document.body.search(\"July 13, 2022\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"March 21, 2023\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The \"Due:\" date is being updated from \"July 13, 2022\" to \"March 21, 2023\".\n# (All other changes in the source diff are auto-generated proofing marks\n# (w:proofErr spell/grammar check tags) and OOXML namespace/schema bumps\n# produced by a newer Word version opening/re-saving the file; they carry\n# no visible text change, so there is nothing else to edit here.)\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n    \"July 13, 2022\",  # FindText\n    $false,            # MatchCase\n    $false,            # MatchWholeWord\n    $false,            # MatchWildcards\n    $false,            # MatchSoundsLike\n    $false,            # MatchAllWordForms\n    $true,             # Forward\n    1,                 # Wrap (wdFindContinue)\n    $false,            # Format\n    \"March 21, 2023\",  # ReplaceWith\n    2                  # Replace (wdReplaceAll)\n)\n"}
